$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty score cells with the value 5
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 5

$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 5

$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 5

$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 5

# Move the scroll position / active selection down to row 28 (as in the
# author's session) and select F28, matching the updated sheetView state.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F28").Select()
